$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '35.321.30'
Set-TextValue 'E2' '  +0.41%  '
Set-TextValue 'D3' '1.885.70'
Set-TextValue 'E3' '  -0.77%  '
Set-TextValue 'E4' '  -0.65%  '
Set-TextValue 'D5' '246.02'
Set-TextValue 'E5' '  -2.55%  '
Set-TextValue 'D6' '0.690'
Set-TextValue 'E6' '  -1.38%  '
Set-TextValue 'E7' '  -0.74%  '
Set-TextValue 'D8' '43.14'
Set-TextValue 'E8' '  +5.43%  '
Set-TextValue 'E9' '  -1.84%  '
Set-TextValue 'D10' '53.67'
Set-TextValue 'E10' '  +1.51%  '
Set-TextValue 'D11' '0.0740'
Set-TextValue 'E11' '  -1.61%  '
Set-TextValue 'D12' '0.0971'
Set-TextValue 'E12' '  -1.14%  '
Set-TextValue 'D13' '13.23'
Set-TextValue 'E13' '  +1.48%  '
Set-TextValue 'D14' '2.162.13'
Set-TextValue 'E14' '  -0.58%  '
Set-TextValue 'D15' '0.754'
Set-TextValue 'E15' '  +2.56%  '
Set-TextValue 'B16' 'WrappedEther'
Set-TextValue 'C16' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D16' '1.904.18'
Set-TextValue 'E16' '  -0.28%  '
Set-TextValue 'B17' 'Polkadot'
Set-TextValue 'C17' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D17' '4.88'
Set-TextValue 'E17' '  -1.39%  '
Set-TextValue 'D18' '35.341.58'
Set-TextValue 'E18' '  +0.47%  '
Set-TextValue 'D19' '72.98'
Set-TextValue 'E19' '  -1.01%  '
Set-TextValue 'D20' '0.0₃0821'
Set-TextValue 'E20' '  -1.45%  '
Set-TextValue 'D21' '244.11'
Set-TextValue 'E21' '  +0.83%  '
Set-TextValue 'D22' '12.75'
Set-TextValue 'E22' '  -1.39%  '
Set-TextValue 'D23' '4.93'
Set-TextValue 'E23' '  -2.08%  '
Set-TextValue 'D24' '2.68'
Set-TextValue 'E24' '  +11.45%  '
Set-TextValue 'E25' '  -0.79%  '
Set-TextValue 'D26' '2.13'
Set-TextValue 'E26' '  -5.53%  '
Set-TextValue 'D27' '166.04'
Set-TextValue 'E27' '  -0.32%  '
Set-TextValue 'D28' '8.48'
Set-TextValue 'E28' '  -1.17%  '
Set-TextValue 'D29' '18.27'
Set-TextValue 'E29' '  -1.12%  '
Set-TextValue 'E30' '  -2.16%  '
Set-TextValue 'D31' '4.128.46'
Set-TextValue 'E32' '  +11.75%  '
Set-TextValue 'E33' '  -1.52%  '
Set-TextValue 'D34' '0.0581'
Set-TextValue 'E34' '  -3.48%  '
Set-TextValue 'D35' '4.16'
Set-TextValue 'E35' '  -0.90%  '
Set-TextValue 'B36' 'BinanceUSD'
Set-TextValue 'C36' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D36' '1.00'
Set-TextValue 'E36' '  -0.71%  '
Set-TextValue 'B37' 'WEMIXToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D37' '1.83'
Set-TextValue 'E37' '  -11.73%  '
Set-TextValue 'E38' '  -1.26%  '
Set-TextValue 'D39' '1.95'
Set-TextValue 'E39' '  -2.58%  '
Set-TextValue 'D40' '0.0691'
Set-TextValue 'E40' '  +6.84%  '
Set-TextValue 'D41' '0.0219'
Set-TextValue 'E41' '  +2.70%  '
Set-TextValue 'D42' '17.14'
Set-TextValue 'E42' '  +0.00%  '
Set-TextValue 'D43' '96.41'
Set-TextValue 'E43' '  -3.69%  '
Set-TextValue 'D44' '1.07'
Set-TextValue 'E44' '  -2.48%  '
Set-TextValue 'D45' '1.298.51'
Set-TextValue 'E45' '  -1.21%  '
Set-TextValue 'E46' '  -4.48%  '
Set-TextValue 'D47' '0.0795'
Set-TextValue 'E47' '  +7.54%  '
Set-TextValue 'D48' '2.40'
Set-TextValue 'E48' '  -1.07%  '
Set-TextValue 'D49' '12.25'
Set-TextValue 'E49' '  +3.73%  '
Set-TextValue 'D50' '2.73'
Set-TextValue 'E50' '  -0.79%  '
Set-TextValue 'D51' '6.24'
Set-TextValue 'E51' '  -5.33%  '
